$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2778.625
$ws.Range("I34").Value = 2778.625
$ws.Range("K34").Value = 2778.625
$ws.Range("M34").Value = -2575.625
$ws.Range("H36").Value = 2778.625
$ws.Range("I36").Value = 2778.625
$ws.Range("K36").Value = 2778.625
$ws.Range("M36").Value = -2063.625
$ws.Range("H135").Value = 1539.0952
$ws.Range("I135").Value = 1572.925
$ws.Range("K135").Value = 14156.325
$ws.Range("M135").Value = -11621.325
$ws.Range("H138").Value = 3031.9138
$ws.Range("I138").Value = 1221
$ws.Range("J138").Value = 5111.1113
$ws.Range("K138").Value = 3663
$ws.Range("L138").Value = 15333.3339
$ws.Range("M138").Value = 1477
$ws.Range("N138").Value = -25613.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -825
$ws.Range("N14").ClearContents()
$ws.Range("H32").Value = 705.21
$ws.Range("I32").Value = 725.3196
$ws.Range("K32").Value = 725.3196
$ws.Range("M32").Value = -438.3196
$ws.Range("H63").Value = 4722
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 5000
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 4722
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 25000
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 4118.278
$ws.Range("I74").Value = 2782
$ws.Range("J74").Value = 6790.8335
$ws.Range("K74").Value = 2782
$ws.Range("L74").Value = 6790.8335
$ws.Range("M74").Value = -1908
$ws.Range("N74").Value = -8538.833500000001
$ws.Range("H77").Value = 4118.278
$ws.Range("I77").Value = 2782
$ws.Range("J77").Value = 6790.8335
$ws.Range("K77").Value = 13910
$ws.Range("L77").Value = 33954.1675
$ws.Range("M77").Value = -9542
$ws.Range("N77").Value = -42690.1675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 5250
$ws.Range("J14").Value = 5250
$ws.Range("L14").Value = 5250
$ws.Range("N14").Value = -5594
$ws.Range("H18").Value = 9999
$ws.Range("J18").Value = 9999
$ws.Range("L18").Value = 9999
$ws.Range("N18").Value = -11057
$ws.Range("H20").Value = 2798.7778
$ws.Range("I20").Value = 1242.875
$ws.Range("J20").Value = 4043.5
$ws.Range("K20").Value = 1242.875
$ws.Range("L20").Value = 4043.5
$ws.Range("M20").Value = -995.875
$ws.Range("N20").Value = -4537.5
$ws.Range("H94").Value = 15585.348
$ws.Range("I94").Value = 17538.422
$ws.Range("K94").Value = 17538.422
$ws.Range("M94").Value = -17087.422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2768.3333
$ws.Range("I11").Value = 3805
$ws.Range("J11").Value = 2250
$ws.Range("K11").Value = 3805
$ws.Range("L11").Value = 2250
$ws.Range("M11").Value = -3665
$ws.Range("N11").Value = -2530
$ws.Range("H19").Value = 831.8
$ws.Range("I19").Value = 831.8
$ws.Range("K19").Value = 831.8
$ws.Range("M19").Value = -661.8
$ws.Range("H22").Value = 1182.4584
$ws.Range("I22").Value = 991.5
$ws.Range("J22").Value = 1449.8
$ws.Range("K22").Value = 991.5
$ws.Range("L22").Value = 1449.8
$ws.Range("M22").Value = -641.5
$ws.Range("N22").Value = -2149.8
$ws.Range("H24").Value = 831.8
$ws.Range("I24").Value = 831.8
$ws.Range("K24").Value = 831.8
$ws.Range("M24").Value = -661.8
$ws.Range("H134").Value = 3450.913
$ws.Range("I134").Value = 1010.1724
$ws.Range("J134").Value = 7614.5293
$ws.Range("K134").Value = 3030.5172
$ws.Range("L134").Value = 22843.5879
$ws.Range("M134").Value = -495.5172000000002
$ws.Range("N134").Value = -27913.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5414.6
$ws.Range("I56").Value = 5414.6
$ws.Range("K56").Value = 5414.6
$ws.Range("M56").Value = -4884.6
$ws.Range("H113").Value = 866.3714
$ws.Range("I113").Value = 554.5
$ws.Range("J113").Value = 991.12
$ws.Range("K113").Value = 1663.5
$ws.Range("L113").Value = 2973.36
$ws.Range("M113").Value = 506.5
$ws.Range("N113").Value = -7313.360000000001
$ws.Range("H140").Value = 340612.97
$ws.Range("I140").Value = 364713.06
$ws.Range("J140").Value = 3211.5
$ws.Range("K140").Value = 1094139.18
$ws.Range("L140").Value = 9634.5
$ws.Range("M140").Value = -1088959.18
$ws.Range("N140").Value = -19994.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 8053.722
$ws.Range("I17").Value = 799.5
$ws.Range("J17").Value = 8960.5
$ws.Range("K17").Value = 799.5
$ws.Range("L17").Value = 8960.5
$ws.Range("M17").Value = -631.5
$ws.Range("N17").Value = -9296.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8870.975
$ws.Range("I22").Value = 12686.667
$ws.Range("K22").Value = 12686.667
$ws.Range("M22").Value = -12391.667
$ws.Range("H27").Value = 8870.975
$ws.Range("I27").Value = 12686.667
$ws.Range("K27").Value = 12686.667
$ws.Range("M27").Value = -12579.667
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H61").Value = 4444.185
$ws.Range("I61").Value = 2444.8333
$ws.Range("K61").Value = 2444.8333
$ws.Range("M61").Value = -2242.8333
$ws.Range("H68").Value = 4668.875
$ws.Range("I68").Value = 4071.2856
$ws.Range("J68").Value = 5133.6665
$ws.Range("K68").Value = 4071.2856
$ws.Range("L68").Value = 5133.6665
$ws.Range("M68").Value = -3322.2856
$ws.Range("N68").Value = -6631.6665
$ws.Range("H71").Value = 4668.875
$ws.Range("I71").Value = 4071.2856
$ws.Range("J71").Value = 5133.6665
$ws.Range("K71").Value = 20356.428
$ws.Range("L71").Value = 25668.3325
$ws.Range("M71").Value = -16612.428
$ws.Range("N71").Value = -33156.3325
$ws.Range("H93").Value = 5586.4546
$ws.Range("I93").Value = 6524.8887
$ws.Range("K93").Value = 6524.8887
$ws.Range("M93").Value = -5276.8887
$ws.Range("H113").Value = 4444.185
$ws.Range("I113").Value = 2444.8333
$ws.Range("K113").Value = 2444.8333
$ws.Range("M113").Value = -274.8332999999998
$ws.Range("H122").Value = 6887.263
$ws.Range("I122").Value = 5718.697
$ws.Range("K122").Value = 17156.091
$ws.Range("M122").Value = -14706.091
$ws.Range("H132").Value = 577351.1
$ws.Range("I132").Value = 1065725.2
$ws.Range("K132").Value = 3197175.6
$ws.Range("M132").Value = -3194645.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 276153.44
$ws.Range("J62").Value = 4225
$ws.Range("L62").Value = 4225
$ws.Range("N62").Value = -5473
$ws.Range("H65").Value = 276153.44
$ws.Range("J65").Value = 4225
$ws.Range("L65").Value = 21125
$ws.Range("N65").Value = -27365
